$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 194 (existing data from row 194 down shifts to 195..278)
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,N,O,Q are identical to the surrounding "Ajo / Chino / Primera"
# records for this market, so only the date + volume/price columns change.
$ws.Cells.Item(194, 1).Value = 7
$ws.Cells.Item(194, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(194, 3).Value = "Ñuble"
$ws.Cells.Item(194, 4).Value = 44837
$ws.Cells.Item(194, 5).Value = 16
$ws.Cells.Item(194, 6).Value = 100112003
$ws.Cells.Item(194, 7).Value = "Ajo"
$ws.Cells.Item(194, 8).Value = "Chino"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 60
$ws.Cells.Item(194, 11).Value = 20000
$ws.Cells.Item(194, 12).Value = 21000
$ws.Cells.Item(194, 13).Value = 20500
$ws.Cells.Item(194, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(194, 15).Value = "China"
$ws.Cells.Item(194, 16).Value = 2050
$ws.Cells.Item(194, 17).Value = 10
$ws.Cells.Item(194, 18).Value = "Hortaliza"
